# Apply symbol-list refresh (prices, 1h volume %, and hour column) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume(1h)) updates, written as literal text (apostrophe-prefixed)
# to preserve the inline-string cell type used throughout this sheet.
$updates = @(
    @{ Cell = "D2"; Value = "'298.53" },
    @{ Cell = "E2"; Value = "'-6.21%" },
    @{ Cell = "D3"; Value = "'35.03" },
    @{ Cell = "E3"; Value = "'-3.21%" },
    @{ Cell = "D4"; Value = "'4.973" },
    @{ Cell = "E4"; Value = "'-2.59%" },
    @{ Cell = "D5"; Value = "'0.07868" },
    @{ Cell = "E5"; Value = "'-2.51%" },
    @{ Cell = "D6"; Value = "'1.900" },
    @{ Cell = "E6"; Value = "'-12.13%" },
    @{ Cell = "E7"; Value = "'-2.47%" },
    @{ Cell = "D8"; Value = "'7.722" },
    @{ Cell = "E8"; Value = "'-4.13%" },
    @{ Cell = "D9"; Value = "'2.908" },
    @{ Cell = "E9"; Value = "'3.89%" },
    @{ Cell = "D10"; Value = "'0.9231" },
    @{ Cell = "E10"; Value = "'-0.60%" },
    @{ Cell = "D11"; Value = "'0.1089" },
    @{ Cell = "E11"; Value = "'7.81%" },
    @{ Cell = "D12"; Value = "'0.1811" },
    @{ Cell = "E12"; Value = "'-3.89%" },
    @{ Cell = "D13"; Value = "'0.09184" },
    @{ Cell = "E13"; Value = "'-0.33%" },
    @{ Cell = "D14"; Value = "'0.03551" },
    @{ Cell = "E14"; Value = "'-1.65%" },
    @{ Cell = "D15"; Value = "'0.09876" },
    @{ Cell = "D16"; Value = "'0.001395" },
    @{ Cell = "E16"; Value = "'-2.89%" },
    @{ Cell = "D17"; Value = "'0.005805" },
    @{ Cell = "E17"; Value = "'2.75%" },
    @{ Cell = "D18"; Value = "'3.486" },
    @{ Cell = "E18"; Value = "'0.70%" },
    @{ Cell = "D19"; Value = "'0.3440" },
    @{ Cell = "E19"; Value = "'1.97%" },
    @{ Cell = "D20"; Value = "'0.1309" },
    @{ Cell = "E20"; Value = "'-2.33%" },
    @{ Cell = "D21"; Value = "'5.070" },
    @{ Cell = "E21"; Value = "'-0.06%" },
    @{ Cell = "D22"; Value = "'0.2400" },
    @{ Cell = "E22"; Value = "'8.99%" },
    @{ Cell = "D23"; Value = "'0.04531" },
    @{ Cell = "E23"; Value = "'-1.53%" },
    @{ Cell = "D24"; Value = "'0.001214" },
    @{ Cell = "E24"; Value = "'-2.29%" },
    @{ Cell = "D25"; Value = "'0.004587" },
    @{ Cell = "E25"; Value = "'-3.40%" },
    @{ Cell = "E26"; Value = "'-3.79%" },
    @{ Cell = "E27"; Value = "'-6.83%" },
    @{ Cell = "E39"; Value = "'-4.48%" },
    @{ Cell = "D40"; Value = "'0.04664" },
    @{ Cell = "E40"; Value = "'-6.13%" },
    @{ Cell = "D41"; Value = "'0.007577" },
    @{ Cell = "E41"; Value = "'-3.13%" },
    @{ Cell = "D42"; Value = "'0.009571" },
    @{ Cell = "E42"; Value = "'22.18%" },
    @{ Cell = "D43"; Value = "'0.1319" },
    @{ Cell = "E43"; Value = "'-5.79%" },
    @{ Cell = "D44"; Value = "'0.002121" },
    @{ Cell = "E44"; Value = "'1.18%" },
    @{ Cell = "D45"; Value = "'0.01084" },
    @{ Cell = "E45"; Value = "'-6.96%" },
    @{ Cell = "D46"; Value = "'0.00005992" },
    @{ Cell = "E46"; Value = "'-4.77%" },
    @{ Cell = "D47"; Value = "'0.00000000750" },
    @{ Cell = "E47"; Value = "'0.05%" },
    @{ Cell = "E48"; Value = "'123.05%" },
    @{ Cell = "E49"; Value = "'-31.36%" },
    @{ Cell = "D50"; Value = "'0.00002101" },
    @{ Cell = "E50"; Value = "'0.05%" },
    @{ Cell = "D51"; Value = "'0.0002001" },
    @{ Cell = "E51"; Value = "'0.05%" }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Column G (Hora) goes from 2 to 3 for every data row (2-51)
for ($r = 2; $r -le 51; $r++) {
    $ws.Range("G$r").Value = "'3"
}
